{"js": "// Fix calculation of referencia_comissao: correct payment-term labels in the\n// commission table from ranges (\"X a Y\") to explicit installment lists\n// (\"X/Y/Z\"), matching the corrected CalculationEngine percentage logic.\nconst replacements = [\n  { find: \"30 a 90\", replace: \"30/60/90\" },\n  { find: \"15 a 45\", replace: \"15/30/45\" },\n  { find: \"30 a 60\", replace: \"30/45/60\" },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Fix calculation of referencia_comissao: correct payment-term labels in the\n# commission table from ranges (\"X a Y\") to explicit installment lists\n# (\"X/Y/Z\"), matching the corrected CalculationEngine percentage logic.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"30 a 90\"; Replace = \"30/60/90\" },\n    @{ Find = \"15 a 45\"; Replace = \"15/30/45\" },\n    @{ Find = \"30 a 60\"; Replace = \"30/45/60\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Find\n    $find.Replacement.Text = $pair.Replace\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    # wdReplaceAll = 2 replaces every occurrence in the range in one call\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
